$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" query saved in B2 used to return an extra `Cohort`
# column (sourced from "OPTIONAL MATCH (co:cohort)<-[*]-(c)"). That
# trailing RETURN line - and the now-dangling trailing comma left on
# the preceding "Response to Treatment" line - is removed here, so the
# query lines up with the other tabs (which never returned a Cohort
# column).
$b2 = $ws.Range("B2")
$current = $b2.Value()

$cohortLine = ",`n        coalesce(co.cohort_description, '') AS ``Cohort``"
if ($current.EndsWith($cohortLine)) {
    $updated = $current.Substring(0, $current.Length - $cohortLine.Length)
    $b2.Value = $updated
}

# Move the saved selection/active cell onto the edited cell.
$b2.Select()
